$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- TextBox 125 ("Imputation Approach") ---
$shImputation = $s.Shapes.Item("TextBox 125")
$shImputation.Height = 96.93748031496062
$paraImputation = $shImputation.TextFrame.TextRange.Paragraphs(2)
$paraImputation.Runs(1).Text = "Our study employed mode imputation for categorical  features, serving as a reference for comparison. Additionally, we leveraged the Naïve Bayes imputation technique to impute missing categorical features feature ."

# --- TextBox 134 ("Data Preprocessing") ---
$shPreprocessing = $s.Shapes.Item("TextBox 134")
$paraPreprocessing = $shPreprocessing.TextFrame.TextRange.Paragraphs(2)
$paraPreprocessing.Runs(1).Text = "The data was split into a 30:70 test: train split. Additionally, three copies of the training subset were generated and induced with missing values at varying proportions (10%, 40%, and 70%). "

# --- TextBox 153 ("Data Exploration") : split into 3 runs ---
$shExploration = $s.Shapes.Item("TextBox 153")
$paraExploration = $shExploration.TextFrame.TextRange.Paragraphs(2)
$subExploration = $paraExploration.Characters(145, 26)
$subExploration.Text = "skewed distribution "

# --- TextBox 160 ("Visualization") ---
$shViz = $s.Shapes.Item("TextBox 160")
$paraViz = $shViz.TextFrame.TextRange.Paragraphs(2)
$paraViz.Runs(1).Text = "Line plots and bar graphs were crafted to vividly illustrate the performance disparity between the two models when trained with data imputed using distinct methods."

# --- TextBox 165 ("Model Comparison") ---
$shModel = $s.Shapes.Item("TextBox 165")
$shModel.Height = 87.24377952755906
$paraModel = $shModel.TextFrame.TextRange.Paragraphs(2)
$paraModel.Runs(1).Text = "The results from the K-NN and decision tree classifiers are evaluated using performance metrices including accuracy, precision, recall,  F1-score."
